$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H, matching the style of the other header cells (e.g. G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# "Save" flag values for rows 2-21 (1 when sum == max sum value, else 0)
$saveValues = @(0, 1, 0, 0, 1, 0, 0, 1, 0, 0, 1, 0, 0, 0, 0, 0, 1, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
